# Generate Report for Handback
# Update timestamp values on each sheet to reflect the new handback report generation time.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: G2 = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-22 13:08:23"

# "zh-cn" sheet: H2 = "Correspond Handoff Datetime", K2 = "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-22 13:08:18"
$wsZhCn.Range("K2").Value = "2016-08-22 13:08:35"

# "de-de" sheet: K2 = "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-22 13:08:42"
